$wb = $excel.ActiveWorkbook

# --- Nets Forecast sheet ---
$netsWs = $wb.Worksheets.Item("Nets Forecast")
$netsWs.Range("B2").Value = 10.33695793151855
$netsWs.Range("B3").Value = 27.75533866882324
$netsWs.Range("B4").Value = 16.16401100158691
$netsWs.Range("B5").Value = 2.293780088424683
$netsWs.Range("B6").Value = 1.674256324768066
$netsWs.Range("B7").Value = 5.555194854736328
$netsWs.Range("B8").Value = 19.31033897399902
$netsWs.Range("B9").Value = 8.197628974914551
$netsWs.Range("B10").Value = 9.649618148803711
$netsWs.Range("B11").Value = 16.25031661987305
$netsWs.Range("B12").Value = 117.1874389648438

# --- Cavaliers Forecast sheet ---
$cavsWs = $wb.Worksheets.Item("Cavaliers Forecast")
$cavsWs.Range("B2").Value = 15.86654758453369
$cavsWs.Range("B3").Value = 15.42139530181885
$cavsWs.Range("B4").Value = 8.723228454589844
$cavsWs.Range("B5").Value = 23.49494934082031
$cavsWs.Range("B6").Value = 5.346555709838867
$cavsWs.Range("B7").Value = 8.939894676208496
$cavsWs.Range("B8").Value = 12.07424163818359
$cavsWs.Range("B9").Value = 6.25080394744873
$cavsWs.Range("B10").Value = 24.57846260070801
$cavsWs.Range("B11").Value = 120.6960754394531
